$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all current cell contents (keeps per-cell styles, e.g. the bold/border
# header style on row 1). The shared-strings table is rebuilt on save in the
# exact order strings are (re)written below, so we control that order here.
$ws.Cells.ClearContents()

$headers = @(
    "Sending cluster","Ligand symbol","Receptor symbol","Target cluster",
    "Ligand-expressing cells","Ligand detection rate","Ligand average expression value",
    "Ligand total expression value","Ligand derived specificity of average expression value",
    "Ligand derived specificity of total expression value","Receptor-expressing cells",
    "Receptor detection rate","Receptor average expression value","Receptor total expression value",
    "Receptor derived specificity of average expression value","Receptor derived specificity of total expression value",
    "Edge average expression weight","Edge total expression weight",
    "Edge average expression derived specificity","Edge total expression derived specificity"
)
for ($c = 0; $c -lt $headers.Count; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows 2-7. Columns A-D are text (shared strings), E-T are numeric.
# Values are written column-by-column so that newly introduced shared strings
# ("sCs" then "ECs") end up ordered exactly as produced by the source tool.
$data = @(
    @("FAPs","Fgf17","Fgfr3","ECs",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,2,0.6666666666666666,2.619953333333334,7.85986,0.6100029538328192,0.6100029538328192,3.748471158815556,33.73624042933999,0.4593887222058784,0.4593887222058785),
    @("FAPs","Fgf17","Fgfr3","FAPs",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,3,1,0.5698483333333333,1.709545,0.1326776176306101,0.1326776176306101,0.8153046144838887,7.337741530354998,0.09991853456721218,0.0999185345672122),
    @("FAPs","Fgf17","Fgfr3","sCs",3,1,1.430739666666667,4.292218999999999,0.7530926191740721,0.7530926191740722,3,1,1.105183,3.315549,0.2573194285365706,0.2573194285365706,1.581229157025666,14.231062413231,0.1937853624009815,0.1937853624009815),
    @("sCs","Fgf17","Fgfr3","ECs",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,2,0.6666666666666666,2.619953333333334,7.85986,0.6100029538328192,0.6100029538328192,1.228965962964444,11.06069366668,0.1506142316269408,0.1506142316269408),
    @("sCs","Fgf17","Fgfr3","FAPs",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,3,1,0.5698483333333333,1.709545,0.1326776176306101,0.1326776176306101,0.2673040763011111,2.40573668671,0.03275908306339787,0.03275908306339788),
    @("sCs","Fgf17","Fgfr3","sCs",2,0.6666666666666666,0.4690793333333333,1.407238,0.2469073808259278,0.2469073808259279,3,1,1.105183,3.315549,0.2573194285365706,0.2573194285365706,0.5184185048513333,4.665766543661999,0.06353406613558917,0.06353406613558918)
)

$colCount = $data[0].Count
for ($c = 0; $c -lt $colCount; $c++) {
    for ($i = 0; $i -lt $data.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, $c + 1).Value = $data[$i][$c]
    }
}
